$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ISIN', 'Stock Name', 'Mutual Fund', 'Jan_2026', 'Dec_2025', 'Nov_2025', 'MoM', 'QoQ'),
    @('INE423A01024', 'Adani Enterprises Limited', 'quant ESG Integration Strategy Fund', 8.711195999999999, 9.043240000000001, 4.153361, -0.3320440000000016, 4.557834999999999),
    @('INE364U01010', 'Adani Green Energy Limited', 'quant ESG Integration Strategy Fund', 7.992836, 8.914787, 0, -0.9219510000000009, 7.992836),
    @('INE202B01038', 'Piramal Finance Ltd', 'quant ESG Integration Strategy Fund', 7.951795, 7.017373, 6.512386, 0.9344219999999996, 1.439408999999999),
    @('INE917I01010', 'Bajaj Auto Limited', 'quant ESG Integration Strategy Fund', 7.515284, 6.851212, 6.550761, 0.664072, 0.9645230000000007),
    @('INE931S01010', 'Adani Energy Solutions Limited', 'quant ESG Integration Strategy Fund', 6.68311, 7.18563, 4.079932, -0.5025199999999996, 2.603178),
    @('INE726G01019', 'ICICI Prudential Life Insurance Co Ltd', 'quant ESG Integration Strategy Fund', 5.858935, 0, 0, 5.858935, 5.858935),
    @('INE090A01021', 'ICICI Bank Limited', 'quant ESG Integration Strategy Fund', 5.626443, 0, 2.057932, 5.626443, 3.568511),
    @('INE406A01037', 'Aurobindo Pharma Limited', 'quant ESG Integration Strategy Fund', 5.354146, 4.911449, 5.014171, 0.4426969999999999, 0.3399749999999999),
    @('INE016A01026', 'Dabur India Limited', 'quant ESG Integration Strategy Fund', 4.123527, 0, 0, 4.123527, 4.123527),
    @('INE237A01036', 'Kotak Mahindra Bank Limited', 'quant ESG Integration Strategy Fund', 4.089286, 0, 0, 4.089286, 4.089286),
    @('INE180C01042', 'Capri Global Capital Limited', 'quant ESG Integration Strategy Fund', 2.77975, 2.708321, 2.743759, 0.07142899999999974, 0.03599100000000011),
    @('INE040A01034', 'HDFC Bank Limited', 'quant ESG Integration Strategy Fund', 2.50462, 0, 0, 2.50462, 2.50462),
    @('INE860A01027', 'HCL Technologies Limited', 'quant ESG Integration Strategy Fund', 0, 1.902312, 0, -1.902312, 0),
    @('INE775A01035', 'Samvardhana Motherson International Ltd', 'quant ESG Integration Strategy Fund', 0, 8.769418999999999, 2.671117, -8.769418999999999, -2.671117),
    @('INE484J01027', 'Godrej Properties Limited', 'quant ESG Integration Strategy Fund', 0, 3.130396, 0.589028, -3.130396, -0.589028),
    @('INE768C01028', 'Zydus Wellness Ltd', 'quant ESG Integration Strategy Fund', 0, 0, 6.717629, 0, -6.717629),
    @('INE669C01036', 'Tech Mahindra Limited', 'quant ESG Integration Strategy Fund', 0, 3.371679, 0, -3.371679, 0),
    @('INE002A01018', 'Reliance Industries Limited', 'quant ESG Integration Strategy Fund', 0, 0, 4.925525, 0, -4.925525),
    @('INE423A20016', 'Adani Enterprises Limited Rights', 'quant ESG Integration Strategy Fund', 0, 0, 0.09669, 0, -0.09669),
    @('INE245A01021', 'Tata Power Company Limited', 'quant ESG Integration Strategy Fund', 0, 0, 6.022069, 0, -6.022069),
    @('INE237A01028', 'Kotak Mahindra Bank Limited', 'quant ESG Integration Strategy Fund', 0, 4.131909, 3.926298, -4.131909, -3.926298),
    @('INE192A01025', 'Tata Consumer Products Ltd', 'quant ESG Integration Strategy Fund', 0, 0, 2.121902, 0, -2.121902),
    @('INE14LE01019', 'Aditya Birla Lifestyle Brands Limited', 'quant ESG Integration Strategy Fund', 0, 0, 3.571092, 0, -3.571092),
    @('INE019A01038', 'JSW Steel Limited', 'quant ESG Integration Strategy Fund', 0, 0, 0.83808, 0, -0.83808),
    @('INE271C01023', 'DLF Limited', 'quant ESG Integration Strategy Fund', 0, 7.00138, 7.256188, -7.00138, -7.256188)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}